$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pair table recomputed with the 'M2' macrophage cluster added
# alongside ECs / FAPs / sCs, so sending x target clusters is now a 4x4
# cross join (16 data rows, A2:T17) instead of the previous 3x4 (12 rows).
$numRows = 16
$numCols = 20
$data = New-Object 'object[,]' $numRows,$numCols

# Row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Tnc"
$data[0,2] = "Itgb1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 8.948174
$data[0,7] = 26.844522
$data[0,8] = 0.0695931738232498
$data[0,9] = 0.0695931738232498
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 153.5290173333333
$data[0,13] = 460.587052
$data[0,14] = 0.3172206968818489
$data[0,15] = 0.317220696881849
$data[0,16] = 1373.804361147683
$data[0,17] = 12364.23925032914
$data[0,18] = 0.02207639509843095
$data[0,19] = 0.02207639509843095

# Row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Tnc"
$data[1,2] = "Itgb1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 8.948174
$data[1,7] = 26.844522
$data[1,8] = 0.0695931738232498
$data[1,9] = 0.0695931738232498
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 168.7997026666667
$data[1,13] = 506.3991080000001
$data[1,14] = 0.3487728915577651
$data[1,15] = 0.3487728915577651
$data[1,16] = 1510.449110609597
$data[1,17] = 13594.04199548638
$data[1,18] = 0.024272212467017
$data[1,19] = 0.024272212467017

# Row 4: ECs -> M2
$data[2,0] = "ECs"
$data[2,1] = "Tnc"
$data[2,2] = "Itgb1"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 8.948174
$data[2,7] = 26.844522
$data[2,8] = 0.0695931738232498
$data[2,9] = 0.0695931738232498
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 68.09032333333333
$data[2,13] = 204.27097
$data[2,14] = 0.1406878008722904
$data[2,15] = 0.1406878008722904
$data[2,16] = 609.2840609029266
$data[2,17] = 5483.556548126339
$data[2,18] = 0.009790910580916058
$data[2,19] = 0.00979091058091606

# Row 5: ECs -> sCs
$data[3,0] = "ECs"
$data[3,1] = "Tnc"
$data[3,2] = "Itgb1"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 8.948174
$data[3,7] = 26.844522
$data[3,8] = 0.0695931738232498
$data[3,9] = 0.0695931738232498
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 93.562673
$data[3,13] = 280.688019
$data[3,14] = 0.1933186106880956
$data[3,15] = 0.1933186106880956
$data[3,16] = 837.215077909102
$data[3,17] = 7534.935701181917
$data[3,18] = 0.0134536556768858
$data[3,19] = 0.0134536556768858

# Row 6: FAPs -> ECs
$data[4,0] = "FAPs"
$data[4,1] = "Tnc"
$data[4,2] = "Itgb1"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 105.8801323333333
$data[4,7] = 317.640397
$data[4,8] = 0.8234679448457706
$data[4,9] = 0.8234679448457706
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 153.5290173333333
$data[4,13] = 460.587052
$data[4,14] = 0.3172206968818489
$data[4,15] = 0.317220696881849
$data[4,16] = 16255.67267225996
$data[4,17] = 146301.0540503396
$data[4,18] = 0.2612210753238393
$data[4,19] = 0.2612210753238393

# Row 7: FAPs -> FAPs
$data[5,0] = "FAPs"
$data[5,1] = "Tnc"
$data[5,2] = "Itgb1"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 105.8801323333333
$data[5,7] = 317.640397
$data[5,8] = 0.8234679448457706
$data[5,9] = 0.8234679448457706
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 168.7997026666667
$data[5,13] = 506.3991080000001
$data[5,14] = 0.3487728915577651
$data[5,15] = 0.3487728915577651
$data[5,16] = 17872.53485617399
$data[5,17] = 160852.8137055659
$data[5,18] = 0.2872032962289896
$data[5,19] = 0.2872032962289896

# Row 8: FAPs -> M2
$data[6,0] = "FAPs"
$data[6,1] = "Tnc"
$data[6,2] = "Itgb1"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 105.8801323333333
$data[6,7] = 317.640397
$data[6,8] = 0.8234679448457706
$data[6,9] = 0.8234679448457706
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 68.09032333333333
$data[6,13] = 204.27097
$data[6,14] = 0.1406878008722904
$data[6,15] = 0.1406878008722904
$data[6,16] = 7209.412445152788
$data[6,17] = 64884.71200637509
$data[6,18] = 0.115851894249176
$data[6,19] = 0.115851894249176

# Row 9: FAPs -> sCs
$data[7,0] = "FAPs"
$data[7,1] = "Tnc"
$data[7,2] = "Itgb1"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 105.8801323333333
$data[7,7] = 317.640397
$data[7,8] = 0.8234679448457706
$data[7,9] = 0.8234679448457706
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 93.562673
$data[7,13] = 280.688019
$data[7,14] = 0.1933186106880956
$data[7,15] = 0.1933186106880956
$data[7,16] = 9906.428198700394
$data[7,17] = 89157.85378830355
$data[7,18] = 0.1591916790437657
$data[7,19] = 0.1591916790437657

# Row 10: M2 -> ECs
$data[8,0] = "M2"
$data[8,1] = "Tnc"
$data[8,2] = "Itgb1"
$data[8,3] = "ECs"
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.1148696666666667
$data[8,7] = 0.344609
$data[8,8] = 0.0008933827928862465
$data[8,9] = 0.0008933827928862465
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 153.5290173333333
$data[8,13] = 460.587052
$data[8,14] = 0.3172206968818489
$data[8,15] = 0.317220696881849
$data[8,16] = 17.63582704474089
$data[8,17] = 158.722443402668
$data[8,18] = 0.0002833995121416276
$data[8,19] = 0.0002833995121416277

# Row 11: M2 -> FAPs
$data[9,0] = "M2"
$data[9,1] = "Tnc"
$data[9,2] = "Itgb1"
$data[9,3] = "FAPs"
$data[9,4] = 2
$data[9,5] = 0.6666666666666666
$data[9,6] = 0.1148696666666667
$data[9,7] = 0.344609
$data[9,8] = 0.0008933827928862465
$data[9,9] = 0.0008933827928862465
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 168.7997026666667
$data[9,13] = 506.3991080000001
$data[9,14] = 0.3487728915577651
$data[9,15] = 0.3487728915577651
$data[9,16] = 19.38996557875245
$data[9,17] = 174.509690208772
$data[9,18] = 0.0003115876999428882
$data[9,19] = 0.0003115876999428882

# Row 12: M2 -> M2
$data[10,0] = "M2"
$data[10,1] = "Tnc"
$data[10,2] = "Itgb1"
$data[10,3] = "M2"
$data[10,4] = 2
$data[10,5] = 0.6666666666666666
$data[10,6] = 0.1148696666666667
$data[10,7] = 0.344609
$data[10,8] = 0.0008933827928862465
$data[10,9] = 0.0008933827928862465
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 68.09032333333333
$data[10,13] = 204.27097
$data[10,14] = 0.1406878008722904
$data[10,15] = 0.1406878008722904
$data[10,16] = 7.821512744525555
$data[10,17] = 70.39361470073
$data[10,18] = 0.0001256880604683109
$data[10,19] = 0.0001256880604683109

# Row 13: M2 -> sCs
$data[11,0] = "M2"
$data[11,1] = "Tnc"
$data[11,2] = "Itgb1"
$data[11,3] = "sCs"
$data[11,4] = 2
$data[11,5] = 0.6666666666666666
$data[11,6] = 0.1148696666666667
$data[11,7] = 0.344609
$data[11,8] = 0.0008933827928862465
$data[11,9] = 0.0008933827928862465
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 93.562673
$data[11,13] = 280.688019
$data[11,14] = 0.1933186106880956
$data[11,15] = 0.1933186106880956
$data[11,16] = 10.74751305995233
$data[11,17] = 96.727617539571
$data[11,18] = 0.0001727075203334199
$data[11,19] = 0.0001727075203334199

# Row 14: sCs -> ECs
$data[12,0] = "sCs"
$data[12,1] = "Tnc"
$data[12,2] = "Itgb1"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 13.635153
$data[12,7] = 40.905459
$data[12,8] = 0.1060454985380935
$data[12,9] = 0.1060454985380935
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 153.5290173333333
$data[12,13] = 460.587052
$data[12,14] = 0.3172206968818489
$data[12,15] = 0.317220696881849
$data[12,16] = 2093.391641279652
$data[12,17] = 18840.52477151687
$data[12,18] = 0.0336398269474371
$data[12,19] = 0.0336398269474371

# Row 15: sCs -> FAPs
$data[13,0] = "sCs"
$data[13,1] = "Tnc"
$data[13,2] = "Itgb1"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 13.635153
$data[13,7] = 40.905459
$data[13,8] = 0.1060454985380935
$data[13,9] = 0.1060454985380935
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 168.7997026666667
$data[13,13] = 506.3991080000001
$data[13,14] = 0.3487728915577651
$data[13,15] = 0.3487728915577651
$data[13,16] = 2301.609772214509
$data[13,17] = 20714.48794993058
$data[13,18] = 0.03698579516181561
$data[13,19] = 0.0369857951618156

# Row 16: sCs -> M2
$data[14,0] = "sCs"
$data[14,1] = "Tnc"
$data[14,2] = "Itgb1"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 13.635153
$data[14,7] = 40.905459
$data[14,8] = 0.1060454985380935
$data[14,9] = 0.1060454985380935
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 68.09032333333333
$data[14,13] = 204.27097
$data[14,14] = 0.1406878008722904
$data[14,15] = 0.1406878008722904
$data[14,16] = 928.42197646947
$data[14,17] = 8355.797788225231
$data[14,18] = 0.01491930798173005
$data[14,19] = 0.01491930798173006

# Row 17: sCs -> sCs
$data[15,0] = "sCs"
$data[15,1] = "Tnc"
$data[15,2] = "Itgb1"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 13.635153
$data[15,7] = 40.905459
$data[15,8] = 0.1060454985380935
$data[15,9] = 0.1060454985380935
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 93.562673
$data[15,13] = 280.688019
$data[15,14] = 0.1933186106880956
$data[15,15] = 0.1933186106880956
$data[15,16] = 1275.741361443969
$data[15,17] = 11481.67225299572
$data[15,18] = 0.02050056844711071
$data[15,19] = 0.02050056844711071

$ws.Range("A2:T17").Value = $data
